$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 598. This shifts the former rows 598-600
# (old weekly entries) down to rows 603-605, preserving their content and
# formatting, while leaving blank rows 598-602 for the new week's data.
$ws.Range("A598:A602").EntireRow.Insert()

function Set-DataRow {
    param(
        [int]$r,
        [int]$a,
        [string]$b,
        [string]$c,
        [double]$d,
        [int]$e,
        [double]$f,
        [string]$g,
        [string]$h,
        [string]$i,
        [double]$j,
        [double]$k,
        [double]$l,
        [double]$m,
        [string]$n,
        [string]$o,
        [double]$p,
        [double]$q,
        [string]$rr
    )
    $ws.Cells.Item($r,1).Value = $a
    $ws.Cells.Item($r,2).Value = $b
    $ws.Cells.Item($r,3).Value = $c
    $ws.Cells.Item($r,4).Value = $d
    $ws.Cells.Item($r,5).Value = $e
    $ws.Cells.Item($r,6).Value = $f
    $ws.Cells.Item($r,7).Value = $g
    $ws.Cells.Item($r,8).Value = $h
    $ws.Cells.Item($r,9).Value = $i
    $ws.Cells.Item($r,10).Value = $j
    $ws.Cells.Item($r,11).Value = $k
    $ws.Cells.Item($r,12).Value = $l
    $ws.Cells.Item($r,13).Value = $m
    $ws.Cells.Item($r,14).Value = $n
    $ws.Cells.Item($r,15).Value = $o
    $ws.Cells.Item($r,16).Value = $p
    $ws.Cells.Item($r,17).Value = $q
    $ws.Cells.Item($r,18).Value = $rr
}

# New row 598 : Ají, Americana (o), Primera
Set-DataRow 598 10 "Vega Modelo de Temuco" "La Araucanía" 44595 9 100112021 "Ají" "Americana (o)" "Primera" 100 17000 17000 17000 "$/caja 15 kilos" "Región del Maule" 1133 15 "Hortaliza"

# New row 599 : Ají, Americana (o), Primera
Set-DataRow 599 10 "Vega Modelo de Temuco" "La Araucanía" 44595 9 100112021 "Ají" "Americana (o)" "Primera" 100 20000 20000 20000 "$/caja 25 kilos" "Región del Maule" 800 25 "Hortaliza"

# New row 600 : Ají, Chilena(o), Primera
Set-DataRow 600 10 "Vega Modelo de Temuco" "La Araucanía" 44595 9 100112021 "Ají" "Chilena(o)" "Primera" 50 35000 35000 35000 "$/saco 25 kilos" "Región del Maule" 1400 25 "Hortaliza"

# New row 601 : Ají, Chilena(o), Segunda
Set-DataRow 601 10 "Vega Modelo de Temuco" "La Araucanía" 44595 9 100112021 "Ají" "Chilena(o)" "Segunda" 20 18000 18000 18000 "$/saco 25 kilos" "Región del Maule" 720 25 "Hortaliza"

# New row 602 : Ají, Inferno, Primera
Set-DataRow 602 10 "Vega Modelo de Temuco" "La Araucanía" 44595 9 100112021 "Ají" "Inferno" "Primera" 150 17000 17000 17000 "$/caja 15 kilos" "Región del Maule" 1133 15 "Hortaliza"
